$wb = $excel.ActiveWorkbook

$black = $wb.Worksheets.Item("Black")
$white = $wb.Worksheets.Item("White")

# Black sheet - batting row (row 2)
$black.Range("B2").Value = 1.0
$black.Range("C2").Value = 10.0
$black.Range("D2").Value = 25.0
$black.Range("E2").Value = 3.0
$black.Range("F2").Value = 5.0
$black.Range("G2").Value = 1.0
$black.Range("H2").Value = 8.0
$black.Range("I2").Value = 8.0
$black.Range("J2").Value = 6.0
$black.Range("K2").Value = 1.0
$black.Range("L2").Value = 1.0
$black.Range("M2").Value = 0.4

# Black sheet - pitching row (row 4)
$black.Range("B4").Value = 1.0
$black.Range("C4").Value = 1.0
$black.Range("D4").Value = 1.0
$black.Range("G4").Value = 2.0
$black.Range("H4").Value = 5.0
$black.Range("I4").Value = 6.0
$black.Range("J4").Value = 1.0
$black.Range("K4").Value = 4.0
$black.Range("L4").Value = 5.0
$black.Range("M4").Value = 5.0
$black.Range("N4").Value = 3.6

# White sheet - batting row (row 2)
$white.Range("B2").Value = 1.0
$white.Range("C2").Value = 5.0
$white.Range("D2").Value = 20.0
$white.Range("F2").Value = 1.0
$white.Range("H2").Value = 2.0
$white.Range("I2").Value = 6.0
$white.Range("J2").Value = 5.0
$white.Range("K2").Value = 4.0
$white.Range("L2").Value = 1.0
$white.Range("M2").Value = 0.25

# White sheet - pitching row (row 4)
$white.Range("B4").Value = 1.0
$white.Range("C4").Value = 1.0
$white.Range("E4").Value = 1.0
$white.Range("G4").Value = 7.0
$white.Range("H4").Value = 10.0
$white.Range("I4").Value = 8.0
$white.Range("J4").Value = 1.0
$white.Range("K4").Value = 1.0
$white.Range("L4").Value = 6.0
$white.Range("M4").Value = 5.0
$white.Range("N4").Value = 12.6
